$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('AB2').Value = 'maa://21246 (91.41), maa://36684 (95.61), ***maa://22731 (6.25)'
$ws.Range('AF2').Value = 'maa://25251 (91.53), ***maa://21730 (26.32), ***maa://39501 (15.62), **maa://36675 (50.0)'
$ws.Range('AB3').Value = 'maa://24390 (94.74)'
$ws.Range('D3').Value = 'maa://36987 (96.08), maa://40192 (98.15), maa://39849 (88.89)'
$ws.Range('H3').Value = 'maa://21247 (98.58), *maa://22748 (60.0)'
$ws.Range('T3').Value = 'maa://24617 (89.74), **maa://20790 (43.48), ***maa://37170 (16.67), maa://45854 (86.36)'
$ws.Range('D4').Value = 'maa://24632 (94.15), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range('T4').Value = 'maa://32509 (95.76), maa://27295 (85.92), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)'
$ws.Range('L5').Value = '*maa://22757 (78.38)'
$ws.Range('H7').Value = '*maa://22763 (66.67)'
$ws.Range('L7').Value = 'maa://28624 (92.79), maa://24957 (97.73)'
$ws.Range('X7').Value = 'maa://22399 (95.48), *maa://22758 (74.29)'
$ws.Range('A8').Value = '更新日期：2025.02.27 13:19:23'
$ws.Range('P8').Value = 'maa://32931 (83.9), *maa://21916 (62.12), maa://23252 (91.18), maa://37496 (96.88), **maa://22759 (45.45)'
$ws.Range('AB9').Value = 'maa://28711 (86.99), ***maa://22740 (5.66), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (96.3), *maa://45044 (66.67)'
$ws.Range('D9').Value = 'maa://22765 (91.67), *maa://21915 (68.97)'
$ws.Range('H9').Value = ''
$ws.Range('P9').Value = 'maa://22736 (83.17)'
$ws.Range('X9').Value = 'maa://26223 (97.92)'
$ws.Range('AF10').Value = '*maa://25021 (54.74), *maa://22733 (61.11), **maa://22761 (50.0)'
$ws.Range('D10').Value = '***maa://25695 (18.62), ***maa://39951 (14.81), ***maa://34206 (20.0), ***maa://39243 (25.0), *maa://45271 (58.33)'
$ws.Range('AB11').Value = 'maa://29912 (97.18), maa://22516 (88.37), *maa://20794 (52.24)'
$ws.Range('L11').Value = 'maa://21287 (89.32)'
$ws.Range('T11').Value = 'maa://22747 (92.5), maa://22501 (97.73), maa://45521 (84.21)'
$ws.Range('X11').Value = 'maa://36713 (97.77)'
$ws.Range('AB12').Value = 'maa://23669 (95.5), maa://36677 (93.65), maa://39872 (91.67)'
$ws.Range('X12').Value = 'maa://22753 (91.01), *maa://21485 (76.43), maa://37962 (90.24)'
$ws.Range('AF13').Value = '**maa://22737 (34.25), maa://39883 (91.67), *maa://39885 (53.33)'
$ws.Range('D13').Value = 'maa://24999 (92.12), maa://36673 (93.24), maa://25001 (85.71)'
$ws.Range('P13').Value = 'maa://22676 (92.74), *maa://22583 (75.0), *maa://22500 (58.7)'
$ws.Range('L14').Value = 'maa://26245 (96.73), maa://21288 (96.3), maa://39841 (95.28), maa://36682 (97.44)'
$ws.Range('T14').Value = 'maa://22521 (94.34), maa://42751 (100.0)'
$ws.Range('AF15').Value = 'maa://21364 (80.78), *maa://36666 (78.18), *maa://22766 (68.64)'
$ws.Range('D15').Value = '*maa://22743 (78.04), maa://22734 (84.17), *maa://30808 (64.18), **maa://36048 (45.0), maa://45058 (92.86)'
$ws.Range('X15').Value = '*maa://38786 (75.0)'
$ws.Range('AB16').Value = 'maa://26228 (95.96)'
$ws.Range('D16').Value = 'maa://21441 (96.4), maa://36679 (94.34), maa://37650 (97.3)'
$ws.Range('T16').Value = 'maa://22729 (95.03), *maa://28648 (69.12), maa://36674 (82.69)'
$ws.Range('H17').Value = 'maa://22430 (88.78), maa://39599 (86.27)'
$ws.Range('T17').Value = '**maa://42324 (50.0)'
$ws.Range('D18').Value = 'maa://24570 (97.35)'
$ws.Range('L18').Value = 'maa://22466 (90.36), *maa://22732 (50.55)'
$ws.Range('AB19').Value = '*maa://30709 (65.83), *maa://36668 (57.5)'
$ws.Range('D20').Value = 'maa://21432 (90.17), maa://25198 (93.64), *maa://20795 (50.77), maa://36680 (91.18)'
$ws.Range('L20').Value = 'maa://41331 (85.06)'
$ws.Range('AB21').Value = 'maa://21443 (81.12), ***maa://23820 (30.0)'
$ws.Range('AF21').Value = 'maa://22524 (94.17), *maa://22432 (77.03)'
$ws.Range('D21').Value = 'maa://21261 (97.67)'
$ws.Range('H21').Value = 'maa://24372 (97.0)'
$ws.Range('P21').Value = 'maa://24381 (85.0)'
$ws.Range('L22').Value = '*maa://27127 (79.66), *maa://22751 (71.01)'
$ws.Range('AF23').Value = 'maa://31489 (94.12)'
$ws.Range('L23').Value = 'maa://39756 (95.68), maa://39875 (94.44)'
$ws.Range('X23').Value = '*maa://28503 (68.75)'
$ws.Range('AF24').Value = 'maa://22523 (86.0), maa://36672 (80.36), maa://29910 (93.22), **maa://21440 (35.71), *maa://45831 (75.0)'
$ws.Range('D24').Value = '*maa://24368 (78.27), *maa://46650 (60.0)'
$ws.Range('X24').Value = 'maa://29988 (84.65), maa://23504 (93.18), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (78.21), ***maa://22815 (23.08)'
$ws.Range('AB25').Value = 'maa://31215 (87.93), maa://24516 (80.22), maa://26001 (87.5)'
$ws.Range('D25').Value = 'maa://29753 (95.17)'
$ws.Range('T25').Value = 'maa://20109 (92.13), maa://22545 (100.0), *maa://42915 (75.0)'
$ws.Range('AB26').Value = 'maa://42235 (94.12)'
$ws.Range('AF26').Value = 'maa://30511 (80.49), *maa://29760 (60.0)'
$ws.Range('D26').Value = 'maa://41802 (94.44)'
$ws.Range('AF28').Value = 'maa://36660 (92.24), *maa://36701 (66.67)'
$ws.Range('D28').Value = 'maa://24465 (90.93), maa://25725 (83.91)'
$ws.Range('T28').Value = '*maa://29765 (63.86), maa://23263 (95.28)'
$ws.Range('X28').Value = 'maa://39929 (90.7), maa://41749 (90.7), ***maa://39723 (13.89)'
$ws.Range('AF29').Value = '*maa://24080 (68.85), maa://42865 (80.95), ***maa://34960 (8.33)'
$ws.Range('L29').Value = 'maa://28432 (93.45), *maa://28440 (80.0), maa://31400 (98.81), *maa://28650 (71.43)'
$ws.Range('AB30').Value = 'maa://42979 (96.97), maa://45822 (100.0), *maa://45045 (80.0)'
$ws.Range('L31').Value = 'maa://35926 (93.47), maa://36258 (85.47), *maa://43904 (72.73)'
$ws.Range('H32').Value = 'maa://21895 (97.01), maa://36667 (97.62), **maa://20793 (38.78), maa://22760 (100.0)'
$ws.Range('T32').Value = 'maa://42859 (95.93), maa://41108 (88.0), maa://41238 (97.12), maa://45523 (100.0)'
$ws.Range('AF34').Value = '*maa://32650 (72.22)'
$ws.Range('T34').Value = 'maa://24526 (93.08)'
$ws.Range('L37').Value = 'maa://45718 (98.55), *maa://47069 (69.23), maa://45789 (100.0)'
$ws.Range('P37').Value = 'maa://21280 (89.4), *maa://21239 (66.67)'
$ws.Range('AF38').Value = 'maa://36697 (86.38)'
$ws.Range('H39').Value = 'maa://36670 (89.11), maa://25199 (84.82), maa://30434 (91.46), ***maa://25036 (16.0), maa://45059 (82.35), *maa://44165 (66.67)'
$ws.Range('P39').Value = 'maa://24709 (91.61), maa://47093 (100.0)'
$ws.Range('T39').Value = 'maa://45788 (81.25), *maa://45790 (75.0), maa://47079 (94.12)'
$ws.Range('H45').Value = 'maa://21229 (84.82), maa://30807 (95.65), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (85.0)'
$ws.Range('T45').Value = '**maa://39364 (38.71)'
$ws.Range('H46').Value = 'maa://35931 (91.95), maa://43901 (92.59)'
$ws.Range('H53').Value = 'maa://32534 (94.08), **maa://32434 (33.33)'
$ws.Range('H55').Value = 'maa://32532 (92.21)'
$ws.Range('H58').Value = '*maa://37964 (60.53)'
$ws.Range('H59').Value = 'maa://31270 (95.35), maa://27746 (82.3)'
$ws.Range('H60').Value = '*maa://40438 (69.35)'
